$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (7->8, and col 20 8->9)
$colsTo8 = @(2,3,5,6,7,8,9,10,11,12,13,15,16,17,21,22,23,24,26,27,28,29,30,32,34)
foreach ($c in $colsTo8) {
    $ws.Columns.Item($c).ColumnWidth = 7.166666666666666
}
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666

# Update data rows 2-5 with new sensor readings
$ws.Cells.Item(2,1).Value = 45078.50694444445
$ws.Cells.Item(2,2).Value = 7.794
$ws.Cells.Item(2,3).Value = 5.793
$ws.Cells.Item(2,4).Value = 1.946
$ws.Cells.Item(2,5).Value = 16.826
$ws.Cells.Item(2,6).Value = 13.11
$ws.Cells.Item(2,7).Value = 4.852
$ws.Cells.Item(2,8).Value = 16.497
$ws.Cells.Item(2,9).Value = 9.57
$ws.Cells.Item(2,10).Value = 4.646
$ws.Cells.Item(2,11).Value = 5.28
$ws.Cells.Item(2,12).Value = 6.967
$ws.Cells.Item(2,13).Value = 7.742
$ws.Cells.Item(2,14).Value = 2.899
$ws.Cells.Item(2,15).Value = 6.391
$ws.Cells.Item(2,16).Value = 8.101000000000001
$ws.Cells.Item(2,17).Value = 5.893
$ws.Cells.Item(2,18).Value = 0.984
$ws.Cells.Item(2,19).Value = 0.93
$ws.Cells.Item(2,20).Value = 89.866
$ws.Cells.Item(2,21).Value = 17.483
$ws.Cells.Item(2,22).Value = 5.899
$ws.Cells.Item(2,23).Value = 10.767
$ws.Cells.Item(2,24).Value = 6.309
$ws.Cells.Item(2,25).Value = 0.916
$ws.Cells.Item(2,26).Value = 10.451
$ws.Cells.Item(2,27).Value = 5.134
$ws.Cells.Item(2,28).Value = 4.948
$ws.Cells.Item(2,29).Value = 6.236
$ws.Cells.Item(2,30).Value = 7.806
$ws.Cells.Item(2,31).Value = 1.833
$ws.Cells.Item(2,32).Value = 14.334
$ws.Cells.Item(2,33).Value = 3.055
$ws.Cells.Item(2,34).Value = 7.294

$ws.Cells.Item(3,1).Value = 45078.51388888889
$ws.Cells.Item(3,2).Value = 22.821
$ws.Cells.Item(3,3).Value = 17.099
$ws.Cells.Item(3,4).Value = 1.568
$ws.Cells.Item(3,5).Value = 49.92
$ws.Cells.Item(3,6).Value = 40.591
$ws.Cells.Item(3,7).Value = 17.319
$ws.Cells.Item(3,8).Value = 65.199
$ws.Cells.Item(3,9).Value = 27.707
$ws.Cells.Item(3,10).Value = 12.843
$ws.Cells.Item(3,11).Value = 17.843
$ws.Cells.Item(3,12).Value = 20.086
$ws.Cells.Item(3,13).Value = 21.466
$ws.Cells.Item(3,14).Value = 6.178
$ws.Cells.Item(3,15).Value = 18.046
$ws.Cells.Item(3,16).Value = 25.278
$ws.Cells.Item(3,17).Value = 15.406
$ws.Cells.Item(3,18).Value = 0.544
$ws.Cells.Item(3,19).Value = 1.025
$ws.Cells.Item(3,20).Value = 267.362
$ws.Cells.Item(3,21).Value = 50.394
$ws.Cells.Item(3,22).Value = 16.657
$ws.Cells.Item(3,23).Value = 33.576
$ws.Cells.Item(3,24).Value = 17.968
$ws.Cells.Item(3,25).Value = 2.389
$ws.Cells.Item(3,26).Value = 33.848
$ws.Cells.Item(3,27).Value = 14.644
$ws.Cells.Item(3,28).Value = 13.144
$ws.Cells.Item(3,29).Value = 15.626
$ws.Cells.Item(3,30).Value = 21.282
$ws.Cells.Item(3,31).Value = 0.915
$ws.Cells.Item(3,32).Value = 59.53
$ws.Cells.Item(3,33).Value = 9.367000000000001
$ws.Cells.Item(3,34).Value = 20.754

$ws.Cells.Item(4,1).Value = 45078.52083333334
$ws.Cells.Item(4,2).Value = 15.678
$ws.Cells.Item(4,3).Value = 11.743
$ws.Cells.Item(4,4).Value = 1.077
$ws.Cells.Item(4,5).Value = 34.374
$ws.Cells.Item(4,6).Value = 27.9
$ws.Cells.Item(4,7).Value = 11.903
$ws.Cells.Item(4,8).Value = 50.108
$ws.Cells.Item(4,9).Value = 19.034
$ws.Cells.Item(4,10).Value = 8.894
$ws.Cells.Item(4,11).Value = 12.215
$ws.Cells.Item(4,12).Value = 13.821
$ws.Cells.Item(4,13).Value = 14.779
$ws.Cells.Item(4,14).Value = 4.234
$ws.Cells.Item(4,15).Value = 12.407
$ws.Cells.Item(4,16).Value = 17.392
$ws.Cells.Item(4,17).Value = 10.621
$ws.Cells.Item(4,18).Value = 0.338
$ws.Cells.Item(4,19).Value = 0.696
$ws.Cells.Item(4,20).Value = 181.533
$ws.Cells.Item(4,21).Value = 34.759
$ws.Cells.Item(4,22).Value = 11.452
$ws.Cells.Item(4,23).Value = 23.151
$ws.Cells.Item(4,24).Value = 12.377
$ws.Cells.Item(4,25).Value = 1.637
$ws.Cells.Item(4,26).Value = 24.918
$ws.Cells.Item(4,27).Value = 10.053
$ws.Cells.Item(4,28).Value = 9.048
$ws.Cells.Item(4,29).Value = 10.728
$ws.Cells.Item(4,30).Value = 14.612
$ws.Cells.Item(4,31).Value = 0.618
$ws.Cells.Item(4,32).Value = 45.809
$ws.Cells.Item(4,33).Value = 6.465
$ws.Cells.Item(4,34).Value = 14.253

$ws.Cells.Item(5,1).Value = 45078.52777777778
$ws.Cells.Item(5,2).Value = 12.83
$ws.Cells.Item(5,3).Value = 9.609999999999999
$ws.Cells.Item(5,4).Value = 0.86
$ws.Cells.Item(5,5).Value = 28.15
$ws.Cells.Item(5,6).Value = 22.86
$ws.Cells.Item(5,7).Value = 9.77
$ws.Cells.Item(5,8).Value = 40.99
$ws.Cells.Item(5,9).Value = 15.58
$ws.Cells.Item(5,10).Value = 7.29
$ws.Cells.Item(5,11).Value = 10.01
$ws.Cells.Item(5,12).Value = 11.31
$ws.Cells.Item(5,13).Value = 12.09
$ws.Cells.Item(5,14).Value = 3.45
$ws.Cells.Item(5,15).Value = 10.15
$ws.Cells.Item(5,16).Value = 14.24
$ws.Cells.Item(5,17).Value = 8.68
$ws.Cells.Item(5,18).Value = 0.25
$ws.Cells.Item(5,19).Value = 0.55
$ws.Cells.Item(5,20).Value = 147.19
$ws.Cells.Item(5,21).Value = 28.43
$ws.Cells.Item(5,22).Value = 9.369999999999999
$ws.Cells.Item(5,23).Value = 18.96
$ws.Cells.Item(5,24).Value = 10.13
$ws.Cells.Item(5,25).Value = 1.34
$ws.Cells.Item(5,26).Value = 20.27
$ws.Cells.Item(5,27).Value = 8.220000000000001
$ws.Cells.Item(5,28).Value = 7.4
$ws.Cells.Item(5,29).Value = 8.75
$ws.Cells.Item(5,30).Value = 11.94
$ws.Cells.Item(5,31).Value = 0.47
$ws.Cells.Item(5,32).Value = 37.37
$ws.Cells.Item(5,33).Value = 5.31
$ws.Cells.Item(5,34).Value = 11.66

# Remove row 6 (dataset now only has 4 data rows)
$ws.Rows.Item(6).Delete()
